$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain value updates (dates/timestamps, temperatures, pressures, etc.)
# These are safe as literal text because the engine does not coerce
# "YYYY-MM-DD HH:MM:SS" or "<number> <unit>" strings into numeric types.
$ws.Range("E2").Value = "2026-02-22 07:18:18"
$ws.Range("E3").Value = "2026-02-22 07:18:20"
$ws.Range("E4").Value = "2026-02-22 07:18:23"
$ws.Range("O4").Value = "6.0 °C"
$ws.Range("E5").Value = "2026-02-22 07:18:25"
$ws.Range("O5").Value = "4.2 °C"
$ws.Range("E6").Value = "2026-02-22 07:18:28"
$ws.Range("E7").Value = "2026-02-22 07:18:30"
$ws.Range("J7").Value = "1027.7 hPa"
$ws.Range("K7").Value = "0.0 MJ/m2"
$ws.Range("O7").Value = "11.6 °C"
$ws.Range("E8").Value = "2026-02-22 07:18:33"
$ws.Range("J8").Value = "1027.5 hPa"
$ws.Range("E9").Value = "2026-02-22 07:18:35"
$ws.Range("E10").Value = "2026-02-22 07:18:38"
$ws.Range("O10").Value = "3.2 °C"
$ws.Range("E11").Value = "2026-02-22 07:18:40"
$ws.Range("E12").Value = "2026-02-22 07:18:42"
$ws.Range("E13").Value = "2026-02-22 07:18:45"
$ws.Range("N13").Value = "-4.2 °C 6:58 TU"
$ws.Range("O13").Value = "-2.8 °C"
$ws.Range("E14").Value = "2026-02-22 07:18:47"
$ws.Range("E15").Value = "2026-02-22 07:18:50"
$ws.Range("E16").Value = "2026-02-22 07:18:52"
$ws.Range("E17").Value = "2026-02-22 07:18:54"
$ws.Range("E18").Value = "2026-02-22 07:18:57"
$ws.Range("E19").Value = "2026-02-22 07:18:59"
$ws.Range("E20").Value = "2026-02-22 07:19:02"
$ws.Range("E21").Value = "2026-02-22 07:19:04"
$ws.Range("E22").Value = "2026-02-22 07:19:07"
$ws.Range("E23").Value = "2026-02-22 07:19:09"
$ws.Range("E24").Value = "2026-02-22 07:19:11"
$ws.Range("E25").Value = "2026-02-22 07:19:14"
$ws.Range("E26").Value = "2026-02-22 07:19:16"
$ws.Range("E27").Value = "2026-02-22 07:19:19"
$ws.Range("E28").Value = "2026-02-22 07:19:21"
$ws.Range("E29").Value = "2026-02-22 07:19:24"
$ws.Range("E30").Value = "2026-02-22 07:19:26"
$ws.Range("E31").Value = "2026-02-22 07:19:28"
$ws.Range("E32").Value = "2026-02-22 07:19:31"
$ws.Range("E33").Value = "2026-02-22 07:19:33"
$ws.Range("E34").Value = "2026-02-22 07:19:36"
$ws.Range("E35").Value = "2026-02-22 07:19:38"
$ws.Range("E36").Value = "2026-02-22 07:19:41"
$ws.Range("E37").Value = "2026-02-22 07:19:43"
$ws.Range("E38").Value = "2026-02-22 07:19:46"
$ws.Range("E39").Value = "2026-02-22 07:19:48"
$ws.Range("E40").Value = "2026-02-22 07:19:50"
$ws.Range("E41").Value = "2026-02-22 07:19:53"
$ws.Range("E42").Value = "2026-02-22 07:19:55"
$ws.Range("E43").Value = "2026-02-22 07:19:57"
$ws.Range("E44").Value = "2026-02-22 07:20:00"
$ws.Range("E45").Value = "2026-02-22 07:20:02"
$ws.Range("E46").Value = "2026-02-22 07:20:05"

# Percentage-looking values need to be forced to text, otherwise Excel
# auto-converts "35%" into the number 0.35 with a percent number format.
# Setting NumberFormat to "@" (Text) first makes the assignment store the
# literal string instead.
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "35%"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "81%"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "64%"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "84%"
